$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217 (shifts rows 217:243 down to 218:244,
# and extends the used range to A1:T244).
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A217").Value = 11
$ws.Range("B217").Value = "Vega Monumental Concepción"
$ws.Range("C217").Value = "Bíobío"
$ws.Range("D217").Value = 45142
$ws.Range("E217").Value = 8
$ws.Range("F217").Value = "Fruta"
$ws.Range("G217").Value = 100102
$ws.Range("H217").Value = "Cítricos"
$ws.Range("I217").Value = 100102004
$ws.Range("J217").Value = "Mandarina"
$ws.Range("K217").Value = "Clementina"
$ws.Range("L217").Value = "Primera"
$ws.Range("M217").Value = 100
$ws.Range("N217").Value = 9000
$ws.Range("O217").Value = 10000
$ws.Range("P217").Value = 9500
$ws.Range("Q217").Value = "$/bandeja 18 kilos"
$ws.Range("R217").Value = "Región de O'Higgins"
$ws.Range("S217").Value = 528
$ws.Range("T217").Value = 18
